$d = $word.ActiveDocument
$r = $d.Content

# Objetivos PT
$old = "Levar aos estudantes conhecimentos básicos sobre a Ciência e Tecnologia dosAlimentos, abrangendo aspectos sobre as instituições envolvidas na produção dealimentos, bem como a habilitação dos profissionais nas respectivas áreas deconcentração da C&T de Alimentos; segurança alimentar / controle de qualidade;legislação vigente; produção de matérias primas, técnicas de processamento; boaspráticas de fabricação; composição dos alimentos, desenvolvimento de novos produtos,características do setor alimentício e sua relação com o meio ambiente. Neste contexto,enfoca-se a evolução do desenvolvimento da C&T dos Alimentos com vistas a atender ademanda por alimentos de qualidade."
$new = "Levar aos estudantes conhecimentos básicos sobre a Ciência e Tecnologia dos^lAlimentos, abrangendo aspectos sobre as instituições envolvidas na produção de^lalimentos, bem como a habilitação dos profissionais nas respectivas áreas de^lconcentração da C&T de Alimentos; segurança alimentar / controle de qualidade;^llegislação vigente; produção de matérias primas, técnicas de processamento; boas^lpráticas de fabricação; composição dos alimentos, desenvolvimento de novos produtos,^lcaracterísticas do setor alimentício e sua relação com o meio ambiente. Neste contexto,^lenfoca-se a evolução do desenvolvimento da C&T dos Alimentos com vistas a atender a^ldemanda por alimentos de qualidade."
$found = $r.Find.Execute($old, $false, $false, $false, $false, $false, $false, 1, $false, $new, 2)
if (-not $found) { Write-Output "MISSING: $old" }

# Objetivos EN
$old = "The aims of this course are focused on the concepts regarding foodscience and technology and relevant aspects in the food segment; considering issuessuch as food demand; food safety, nutritional and fisico chemical characteristics offood. It is also included aspects regarding Food industries. Legislation; Raw Materials;Food Processing Technics; Food Composition, and Development of new foods"
$new = "The aims of this course are focused on the concepts regarding food^lscience and technology and relevant aspects in the food segment; considering issues^lsuch as food demand; food safety, nutritional and fisico chemical characteristics of^lfood. It is also included aspects regarding Food industries. Legislation; Raw Materials;^lFood Processing Technics; Food Composition, and Development of new foods"
$found = $r.Find.Execute($old, $false, $false, $false, $false, $false, $false, 1, $false, $new, 2)
if (-not $found) { Write-Output "MISSING: $old" }

# Programa resumido PT
$old = "Ciência e Tecnologia de Alimentos: conceitos e objetivos (consumo de alimentos,segurança alimentar, integridade e características nutritivas dos alimentos);Generalidades sobre o Setor Alimentício Legislação; Produção de Matérias-Primas;Técnicas de Processamento: Composição dos Alimentos; Lançamentos de Novosprodutos"
$new = "Ciência e Tecnologia de Alimentos: conceitos e objetivos (consumo de alimentos,^lsegurança alimentar, integridade e características nutritivas dos alimentos);^lGeneralidades sobre o Setor Alimentício Legislação; Produção de Matérias-Primas;^lTécnicas de Processamento: Composição dos Alimentos; Lançamentos de Novos^lprodutos"
$found = $r.Find.Execute($old, $false, $false, $false, $false, $false, $false, 1, $false, $new, 2)
if (-not $found) { Write-Output "MISSING: $old" }

# Programa resumido EN
$old = "Food Science and Technology: concepts and objectives (food consumption, food safety, integrity and nutritional characteristics of food);General information about the Food Sector Legislation; Production of Raw Materials;Processing Techniques: Food Composition; New Releases products"
$new = "Food Science and Technology: concepts and objectives (food consumption, food safety, integrity and nutritional characteristics of food);General information about the Food Sector Legislation; Production of Raw Materials;^lProcessing Techniques: Food Composition; New Releases products"
$found = $r.Find.Execute($old, $false, $false, $false, $false, $false, $false, 1, $false, $new, 2)
if (-not $found) { Write-Output "MISSING: $old" }

# Programa PT
$old = "Introdução: conceitos de C&T de alimentos; áreas de concentração e respectivasatribuições dos profissionais; instituições envolvidas com o desenvolvimento da C&TAlimentos; Eembalagens ativas/inteligentes Segurança Alimentar: fatores envolvidos natoxi-infecção alimentar; alimentos orgânicos e Boas Práticas de Fabricação (BPF);microbiologia de alimentos; higiene industrial; análise de perigo dos pontos críticos decontrole - APPCC. Matérias-Primas: importância e características dos segmentosprodutores de matérias primas de origem animal, vegetal, microbiana eaditivos/ingredientes. Processamento de Alimentos: objetivos e caracterização dosdiferentes métodos de processamento dos alimentos abrangendo técnicas deconservação, transformação e melhoria da qualidade; usos e aplicações deaditivos/ingredientes em alimentos. Novos Produtos: fatores que devem serconsiderados no lançamento de novos produtos alimentícios. Composição dosAlimentos: composição e significância dos diferentes compostos encontrados nosalimentos focando em suas propriedades e funções"
$new = "Introdução: conceitos de C&T de alimentos; áreas de concentração e respectivas^latribuições dos profissionais; instituições envolvidas com o desenvolvimento da C&T^lAlimentos; Eembalagens ativas/inteligentes Segurança Alimentar: fatores envolvidos na^ltoxi-infecção alimentar; alimentos orgânicos e Boas Práticas de Fabricação (BPF);^lmicrobiologia de alimentos; higiene industrial; análise de perigo dos pontos críticos de^lcontrole - APPCC. Matérias-Primas: importância e características dos segmentos^lprodutores de matérias primas de origem animal, vegetal, microbiana e^laditivos/ingredientes. Processamento de Alimentos: objetivos e caracterização dos^ldiferentes métodos de processamento dos alimentos abrangendo técnicas de^lconservação, transformação e melhoria da qualidade; usos e aplicações de^laditivos/ingredientes em alimentos. Novos Produtos: fatores que devem ser^lconsiderados no lançamento de novos produtos alimentícios. Composição dos^lAlimentos: composição e significância dos diferentes compostos encontrados nos^lalimentos focando em suas propriedades e funções"
$found = $r.Find.Execute($old, $false, $false, $false, $false, $false, $false, 1, $false, $new, 2)
if (-not $found) { Write-Output "MISSING: $old" }

# Programa EN
$old = "Introduction: concepts of Food Science and Technology; Major fields of Food Scienceand professional skills; Institutions involved with Food & Science TechnologyDevelopment; Food packing; additives and food ingredients.Food Safety: aspects of foodborne diseases organic foods; guidelines for foodproduction; food microbiology; industrial hygiene; food quality control.Raw Material for Food Industries: characteristic of different foodstuff. Animal,vegetal and microbial raw materials as foodstuff. Aditives and food ingredients.Food Processing: objectives and characteristics of different food processing technics.Food preservation, transformations and improvement of food quality methods.Characteristics, use and applications of food additives and ingredients.Novel Food Products: issues involved in a new food products development.Food Composition: food constituents and its significance. Chemical characteristics,properties and functional aspects of different food compounds. Organic, inorganic andtoxic compounds in food composition."
$new = "Introduction: concepts of Food Science and Technology; Major fields of Food Science^land professional skills; Institutions involved with Food & Science Technology^lDevelopment; Food packing; additives and food ingredients.^lFood Safety: aspects of foodborne diseases organic foods; guidelines for food^lproduction; food microbiology; industrial hygiene; food quality control.^lRaw Material for Food Industries: characteristic of different foodstuff. Animal,^lvegetal and microbial raw materials as foodstuff. Aditives and food ingredients.^lFood Processing: objectives and characteristics of different food processing technics.^lFood preservation, transformations and improvement of food quality methods.^lCharacteristics, use and applications of food additives and ingredients.^lNovel Food Products: issues involved in a new food products development.^lFood Composition: food constituents and its significance. Chemical characteristics,^lproperties and functional aspects of different food compounds. Organic, inorganic and^ltoxic compounds in food composition."
$found = $r.Find.Execute($old, $false, $false, $false, $false, $false, $false, 1, $false, $new, 2)
if (-not $found) { Write-Output "MISSING: $old" }

# Norma de recuperacao
$old = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"
$new = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: ^l^lMR = (NF + PR)/2"
$found = $r.Find.Execute($old, $false, $false, $false, $false, $false, $false, 1, $false, $new, 2)
if (-not $found) { Write-Output "MISSING: $old" }

# Bibliografia
$old = "POTTER, N.N., HOTCHKISS, J.H., Food Science - 5ª Ed. Chapman & Hall, 1995.ARAÚJO, J.M.A., Química de Alimentos: Teoria e Prática. 3ª Ed. Editora UFV, 2004.FRAZIER, W.C., Microbiologia de los Alimentos. Editora Acribia, Zaragoza-Epanha,1981. EVANGELISTA, J., Tecnologia de Alimentos, Livraria Atheneu, RJ. 1987.FENEMA, O.R., Principles of Food Science: Part I: Food Chemistry. Marcel Dekker,Inc. NY-USA.1975.BENWART, G.J., Basic Food Microbiology. AVI Publishing Company Inc. USA,1970.TEIXEIRA, L.J.Q; LIMA FILHO, T; SILVA, T; CASSIANO. O. Tecnologia deAlimentos: Processamento Não Térmico Editora: Editora Rubio Ltda. . ISBN:6588340176 ISBN13: 9786588340172/ Páginas: 224 Publicação: Ed.2023ALCARDE, A.R.; D’ARGE, M. R.; SPOLO, M.H.F; Fundamentos de Ciência eTecnologia de Alimentos. E. 2ª - Editora Manole. ISBN: 9788520447147, pags. 480;2019.LAJOLO, F. M & MERCADANTE, A.Z. Química e Bioquímica dos Alimentos -Volume 2 – Ed. 1 - Coleção: Ciência, Tecnologia, Engenharia de Alimentos eNutrição. ISBN-10 ‏ : ‎ 8538808516 ISBN-13 ‏ : ‎ 978-8538808510 . Pag.432. EditoraAtheneu – RJ. 2017."
$new = "POTTER, N.N., HOTCHKISS, J.H., Food Science - 5ª Ed. Chapman & Hall, 1995.^lARAÚJO, J.M.A., Química de Alimentos: Teoria e Prática. 3ª Ed. Editora UFV, 2004.^lFRAZIER, W.C., Microbiologia de los Alimentos. Editora Acribia, Zaragoza-Epanha,^l1981. EVANGELISTA, J., Tecnologia de Alimentos, Livraria Atheneu, RJ. 1987.^lFENEMA, O.R., Principles of Food Science: Part I: Food Chemistry. Marcel Dekker,^lInc. NY-USA.1975.^lBENWART, G.J., Basic Food Microbiology. AVI Publishing Company Inc. USA,1970.^lTEIXEIRA, L.J.Q; LIMA FILHO, T; SILVA, T; CASSIANO. O. Tecnologia de^lAlimentos: Processamento Não Térmico Editora: Editora Rubio Ltda. . ISBN:^l6588340176 ISBN13: 9786588340172/ Páginas: 224 Publicação: Ed.2023^lALCARDE, A.R.; D’ARGE, M. R.; SPOLO, M.H.F; Fundamentos de Ciência e^lTecnologia de Alimentos. E. 2ª - Editora Manole. ISBN: 9788520447147, pags. 480;^l2019.^lLAJOLO, F. M & MERCADANTE, A.Z. Química e Bioquímica dos Alimentos -^lVolume 2 – Ed. 1 - Coleção: Ciência, Tecnologia, Engenharia de Alimentos e^lNutrição. ISBN-10 ‏ : ‎ 8538808516 ISBN-13 ‏ : ‎ 978-8538808510 . Pag.432. Editora^lAtheneu – RJ. 2017."
$found = $r.Find.Execute($old, $false, $false, $false, $false, $false, $false, 1, $false, $new, 2)
if (-not $found) { Write-Output "MISSING: $old" }

Write-Output "DONE"